# Update "LUY_KE_THANG_HE_THONG" report data (report co so)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_HE_THONG")

# Row 7: last_edited_time + several updated numeric properties
$ws.Range("D7").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("T7").Value = 7300000
$ws.Range("W7").Value = 55336000
$ws.Range("AA7").Value = 94894000
$ws.Range("AE7").Value = 150230000
$ws.Range("AH7").Value = 115930000
$ws.Range("AK7").Value = 19
$ws.Range("AN7").Value = 34300000
$ws.Range("AQ7").Value = 123230000

# Rows 8-12: only last_edited_time changes
$ws.Range("D8").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D9").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D10").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D11").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("D12").Value = "2024-07-08T01:58:00.000Z"
